$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 12 (CTTTP04A1.00), shifting rows 12+ down by one.
$ws.Rows(12).Insert()

# The newly inserted row 12 should contain the same NPX values that row 11 (CTTTP03A1.00) has,
# since the CTTTP03A1.00 sample row is being duplicated: one copy keeps its real sample id,
# the other copy (the original row 11 position) becomes the new "non-cimac-control" row.
$ws.Range("A11:CQ11").Copy($ws.Range("A12:CQ12"))

# The copy/paste leaves a stray empty cell in column D (which has no data for this sample);
# clear it so the row shape matches row 11 exactly (no cell in D for these data rows).
$ws.Range("D12").ClearContents()

# Row 11 now represents a non-CIMAC control id instead of the real sample id.
$ws.Range("A11").Value = "non-cimac-control"

# Match the author's final cursor position (now two rows further down than before).
[void]$ws.Range("A15").Select()
